$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.805.77"
$ws.Range("D3").Value = "1.758.27"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.70%  "
$ws.Range("D5").Value = "'237.30"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").Value = "'0.5070"
$ws.Range("E7").Value = "  +3.93%  "
$ws.Range("D8").Value = "'41.24"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("D9").Value = "'0.2656"
$ws.Range("E9").Value = "  +10.38%  "
$ws.Range("D10").Value = "'0.06208"
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("D11").Value = "1.750.63"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").Value = "'0.06940"
$ws.Range("E12").Value = "  +5.09%  "
$ws.Range("D13").Value = "'15.59"
$ws.Range("E13").Value = "  +13.93%  "
$ws.Range("E14").Value = "  +3.69%  "
$ws.Range("D15").Value = "'4.470"
$ws.Range("E15").Value = "  +4.68%  "
$ws.Range("D16").Value = "'77.53"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "25.838.34"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "'0.000006833"
$ws.Range("E20").Value = "  +10.49%  "
$ws.Range("D21").Value = "'11.64"
$ws.Range("E21").Value = "  +8.04%  "
$ws.Range("D22").Value = "1.976.58"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").Value = "'4.069"
$ws.Range("E23").Value = "  +7.07%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'8.162"
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "'5.192"
$ws.Range("E25").Value = "  +2.55%  "
$ws.Range("D26").Value = "'137.87"
$ws.Range("E26").Value = "  +4.28%  "
$ws.Range("E27").Value = "  +3.57%  "
$ws.Range("D28").Value = "'1.822"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("D29").Value = "'15.02"
$ws.Range("E29").Value = "  +6.36%  "
$ws.Range("D30").Value = "'102.76"
$ws.Range("E30").Value = "  +4.17%  "
$ws.Range("D31").Value = "'0.08229"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").Value = "'3.688"
$ws.Range("E32").Value = "  +3.56%  "
$ws.Range("D33").Value = "'3.411"
$ws.Range("E33").Value = "  +9.62%  "
$ws.Range("D34").Value = "'0.04372"
$ws.Range("E34").Value = "  +3.48%  "
$ws.Range("D35").Value = "'0.9997"
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("D36").Value = "'2.655"
$ws.Range("E36").Value = "  +1.39%  "
$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").Value = "'0.6066"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01550"
$ws.Range("E40").Value = "  +8.25%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'1.937"
$ws.Range("E41").Value = "  -6.45%  "
$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").Value = "'103.21"
$ws.Range("E43").Value = "  +2.39%  "
$ws.Range("D44").Value = "'0.3832"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("D45").Value = "'0.7389"
$ws.Range("E45").Value = "  -4.90%  "
$ws.Range("D46").Value = "'4.895"
$ws.Range("E46").Value = "  -5.08%  "
$ws.Range("D47").Value = "'0.05492"
$ws.Range("E47").Value = "  +6.87%  "
$ws.Range("D48").Value = "'0.1081"
$ws.Range("E48").Value = "  +7.13%  "
$ws.Range("D49").Value = "'5.954"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").Value = "'30.05"
$ws.Range("E50").Value = "  +5.06%  "
$ws.Range("D51").Value = "'7.585"
$ws.Range("E51").Value = "  +4.05%  "
